$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Insert a new column before column I (col 9), shifting I:N -> J:O ---
$ws.Columns.Item(9).Insert()

# Give the new column the same default font/style as its neighbours (style 3)
# and roughly match the original column width.
$ws.Columns.Item(9).Font.Name = "Arial"
$ws.Columns.Item(9).Font.Size = 10
$ws.Columns.Item(9).ColumnWidth = 7.3333333

# Header + values for new "Instrument" column
$ws.Cells.Item(1, 9).Value = "Instrument"
for ($r = 2; $r -le 9; $r++) {
    $ws.Cells.Item($r, 9).Value = "Stock"
}

# --- Simplify data validation ranges that got split across insert boundary ---
$gv = $ws.Range("G2:G1048576").Validation
$gType = $gv.Type
$gFormula = $gv.Formula1
$gv.Delete()
$ws.Range("G2:G1048576").Validation.Add(3, 1, 1, '"' + $gFormula + '"')

$mv = $ws.Range("M2:M1048576").Validation
$mFormula = $mv.Formula1
$mv.Delete()
$ws.Range("M2:M1048576").Validation.Add(3, 1, 1, '"' + $mFormula + '"')

$nv = $ws.Range("N2:N1048576").Validation
$nFormula = $nv.Formula1
$nv.Delete()
$ws.Range("N2:N1048576").Validation.Add(3, 1, 1, '"' + $nFormula + '"')

# --- Shift the comments that lived in columns I..N (9..14) one column right ---
for ($col = 14; $col -ge 9; $col--) {
    $srcCell = $ws.Cells.Item(1, $col)
    $cm = $srcCell.Comment
    if ($cm -ne $null) {
        $txt = $cm.Text()
        $cm.Delete()
        $dstCell = $ws.Cells.Item(1, $col + 1)
        $dstCell.AddComment($txt)
    }
}

# --- Selection to match the saved view ---
$ws.Range("I3:I9").Select()
